# Generate Report for Archive
#
# 1. Update status text "Ready for handoff" -> "In Translation" everywhere it appears
#    (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2)
# 2. Narrow the "Status" columns (Overview columns E & F, zh-cn/de-de column C)
#    from width 17.2159881591797 to 13.4101845877511
#    NOTE: Excel's ColumnWidth setter snaps to whole-pixel boundaries (quantized
#    in steps of 1/6 character on this font), so 12.5 is the input value that
#    lands closest to the target stored width of 13.4101845877511.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = $cell.Value2
        if ($oldStatus -eq $v) {
            $cell.Value2 = $newStatus
        }
    }
}

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
